$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (columns B:E)
$ws.Range("B2").Value = 104.4286150598073
$ws.Range("C2").Value = 104.74085481089448
$ws.Range("D2").Value = 106.55847002805434
$ws.Range("E2").Value = 105.43688710742059

# Row 3 data values (columns B:E)
$ws.Range("B3").Value = 104.88524901633632
$ws.Range("C3").Value = 106.0333219377177
$ws.Range("D3").Value = 105.58780782636271
$ws.Range("E3").Value = 105.20251871732125

# Update the selection to reflect the new used region B1:E3
$ws.Range("B1:E3").Select()
